$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.008.85"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.704.55"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4004"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4041"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.473"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.32"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.000"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08824"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.09"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.485"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001353"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.970"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "1.725.22"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.00"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07198"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.78"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.320"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").Value = "24.998.09"
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.406"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.34%  "
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.59"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.069"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +12.65%  "
$ws.Range("E29").Value = "  -2.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "152.29"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.426"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.664"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +20.21%  "
$ws.Range("D33").Value = "1.914.22"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08622"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03169"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.052"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.203"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2920"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09732"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.07"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8261"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.04"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.692"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7386"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.09264"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +13.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.256"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.405"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "139.99"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.01%  "
